$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 8-15 keep their row position but the "name" column shifts because two
# new strings (line7, line8) were inserted into the lookup sequence before
# "extr1". Update B (name), C (from_bus), D (to_bus), E (in_service) for the
# existing rows 8-15, and append two brand-new rows (16, 17).

$rows = @(
    @{ Row = 8;  A = 6;  B = "line7"; C = 14; D = 11; E = $true  },
    @{ Row = 9;  A = 7;  B = "line8"; C = 16; D = 9;  E = $true  },
    @{ Row = 10; A = 8;  B = "extr1"; C = 5;  D = 12; E = $true  },
    @{ Row = 11; A = 9;  B = "extr2"; C = 5;  D = 9;  E = $true  },
    @{ Row = 12; A = 10; B = "extr3"; C = 10; D = 11; E = $false },
    @{ Row = 13; A = 11; B = "extr4"; C = 7;  D = 8;  E = $false },
    @{ Row = 14; A = 12; B = "extr5"; C = 9;  D = 11; E = $false },
    @{ Row = 15; A = 13; B = "extr6"; C = 7;  D = 11; E = $false },
    @{ Row = 16; A = 14; B = "extr7"; C = 5;  D = 7;  E = $false },
    @{ Row = 17; A = 15; B = "extr8"; C = 8;  D = 5;  E = $true  }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    if ($rowNum -gt 15) {
        # Copy the formatting from the last existing data row so the new
        # rows match the established look (bold/bordered index column, etc.)
        $ws.Range("A15").Copy()
        $ws.Range("A$rowNum").PasteSpecial(-4122)
    }

    $ws.Range("A$rowNum").Value = $r.A
    $ws.Range("B$rowNum").Value = $r.B
    $ws.Range("C$rowNum").Value = $r.C
    $ws.Range("D$rowNum").Value = $r.D
    $ws.Range("E$rowNum").Value = $r.E
}
